$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Table/cell left margins: 98 dxa -> 93 dxa (4.9pt -> 4.65pt)
#    affects the table-wide default (tblCellMar) and every row's cell tcMar.
# ---------------------------------------------------------------------------
$tbl = $d.Tables.Item(1)
$tbl.LeftPadding = 4.65

for ($i = 1; $i -le $tbl.Rows.Count; $i++) {
    $cell = $tbl.Cell($i, 1)
    $cell.LeftPadding = 4.65
}

# ---------------------------------------------------------------------------
# 2) Postcondition text: insert "into the database" right after
#    "Imports a whitelist" (before the comma), split across three runs:
#      "Imports a whitelist "
#      "into the database"
#      ", which will be not archived, from the stated filepath into the
#       database, if the file at the stated location is a valid whitelist."
# ---------------------------------------------------------------------------
$oldText = "Imports a whitelist, which will be not archived, from the stated filepath into the database, if the file at the stated location is a valid whitelist."
$newText = "Imports a whitelist into the database, which will be not archived, from the stated filepath into the database, if the file at the stated location is a valid whitelist."

$runA = "Imports a whitelist "
$runB = "into the database"

$found = $d.Content
$found.Find.Execute($oldText)
$start = $found.Start

# Replace the whole sentence first (keeps a single run with the full text).
$found.Text = $newText

# Now split the middle segment ("into the database") into its own run by
# touching an otherwise-inert formatting property on just that sub-range;
# this forces the engine to materialise it as a distinct <w:r>, matching
# the three-run structure introduced by the edit.
$b1 = $start + $runA.Length
$b2 = $b1 + $runB.Length
$segB = $d.Range($b1, $b2)
$segB.Font.Name = "Calibri"
